# Generate Report for Handoff
#
# This script updates the localization-status report to reflect that the
# handoff package is now ready:
#   - Status cells that read "In Translation" become "Ready for handoff"
#   - The associated "Latest HO Xliff Generate Date" / "Latest Handoff
#     Datetime" timestamps are bumped forward a few minutes
#   - The Status columns widen (auto-fit) to accommodate the longer text

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$newStatus = "Ready for handoff"

# --- Overview sheet -------------------------------------------------
$ws1.Range("E2").Value = $newStatus
$ws1.Range("F2").Value = $newStatus
$ws1.Range("G2").Value = "2016-08-31 16:44:43"

# --- zh-cn sheet ------------------------------------------------------
$ws2.Range("C2").Value = $newStatus
$ws2.Range("H2").Value = "2016-08-31 16:44:39"

# --- de-de sheet ------------------------------------------------------
$ws3.Range("C2").Value = $newStatus
$ws3.Range("H2").Value = "2016-08-31 16:44:43"

# --- Widen the Status columns to fit the new, longer text -------------
# (mirrors Excel auto-fitting the column after the text grew)
$newColumnWidth = 16.333333333333336

$ws1.Columns.Item(5).ColumnWidth = $newColumnWidth   # Overview!E (zh-cn status)
$ws1.Columns.Item(6).ColumnWidth = $newColumnWidth   # Overview!F (de-de status)
$ws2.Columns.Item(3).ColumnWidth = $newColumnWidth   # zh-cn!C (Status)
$ws3.Columns.Item(3).ColumnWidth = $newColumnWidth   # de-de!C (Status)
